$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 11 new rows right after the header row (row 1), pushing the existing
# "falling" dataset down by 11 rows.
$ws.Rows("2:12").Insert()
$ws.Range("A2:H12").ClearFormats()

# New sensor-reading rows captured on May 9th, continuing the "falling" series.
$newData = @(
    @(0, "falling", -2.507025241851806, 9.102962493896484, -0.6246470808982849, 0.02334324724790524, 0.01564411001234516, 0.0184190768475939),
    @(100, "falling", -2.812539100646973, 9.111927032470703, -0.3471660315990448, -0.001020592069480415, 0.02263181184123195, -0.0006294894690920053),
    @(200, "falling", -2.887884616851806, 9.205938339233398, -0.2687076330184936, -0.02465064778197097, 0.04819875901065215, -0.01549884358920696),
    @(300, "falling", -2.862497329711914, 9.105484008789062, -0.2620421051979065, 0.01120043709510704, 0.02745168796944902, 0.004529342418763651),
    @(400, "falling", -2.736449241638184, 9.052282333374023, -0.2773746848106384, -0.007870477419800848, 0.04201560953586564, 0.04054804218978412),
    @(500, "falling", -2.651521682739258, 9.127286911010742, 0.0024068877100944, -0.03585853518509281, 0.05588672146564573, 0.04165430539628354),
    @(600, "falling", -2.57433032989502, 9.058347702026367, 0.0778785794973373, 0.01350235557410764, -0.04772198527324475, 0.04825462864303003),
    @(700, "falling", -2.764423370361328, 9.000140190124512, -0.2859586775302887, 0.06278502832098701, -0.1541466276820112, 0.09075818756004657),
    @(800, "falling", -2.818140029907227, 8.811227798461914, -0.6212533712387085, 0.02069492338270688, -0.3678749610738054, 0.1358281277665277),
    @(900, "falling", -3.191051483154297, 9.015185356140137, -0.4772885143756866, -0.1521724931350568, -0.4923760018697599, 0.1619388900878952),
    @(1000, "falling", -2.787490367889404, 8.648155212402344, -1.196338057518005, -0.2075079472326651, -0.5799234204175996, 0.3153703608890859)
)

for ($i = 0; $i -lt $newData.Count; $i++) {
    $r = $i + 2
    $vals = $newData[$i]
    for ($c = 0; $c -lt $vals.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $vals[$c]
    }
}

# The previously-existing rows were pushed down to rows 13:31 by the insert; their
# timestamp column needs to keep counting up in the same 100ms cadence instead of
# restarting at 0, so continue the sequence from the last new row (1000).
for ($r = 13; $r -le 31; $r++) {
    $ws.Cells.Item($r, 1).Value = ($r - 13 + 1) * 100 + 1000
}

# The oldest trailing row (formerly the last row of data) is dropped to keep a
# rolling 30-row window; it is now at row 32 after the insert above.
$ws.Rows("32:32").Delete()

